$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "Actual Time" row for the completed task (row 36) ---
# Copy formatting from the row above (row 35) so the new row's D column
# (date) picks up the existing date-number-format style instead of Excel
# auto-generating a brand new custom numFmt.
$ws.Range("A35:D35").Copy()
$ws.Range("A36:D36").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A36").Value = "Sam"
$ws.Range("B36").Value = "Integrated the login process with the view"
$ws.Range("C36").Value = 3
$ws.Range("D36").Value = "2/27/2025"

# --- Fill in hours-remaining-after (Day 8 / Day 12 / Day 16) for two tasks ---
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 0

$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0

# --- Update the sheet view selection to reflect the new active cell ---
$ws.Range("E36").Select()
